$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F83").Value = 31
$ws.Range("G83").Value = 2070.49

$ws.Range("F92").Value = 122
$ws.Range("G92").Value = 13656.68

$ws.Range("F101").Value = 56
$ws.Range("G101").Value = 1712.48

$ws.Range("B114").Value = 275061.95

$ws.Range("B163").Value = 64329
$ws.Range("E163").Value = 128.32
$ws.Range("F163").Value = 3
$ws.Range("G163").Value = 362.07

$ws.Range("B164").Value = 57552
$ws.Range("E164").Value = 136.86
$ws.Range("F164").Value = -5
$ws.Range("G164").Value = -603.45

$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0

$ws.Range("B222").Value = 56733.17

$ws.Range("F234").Value = 13
$ws.Range("G234").Value = 3844.75

$ws.Range("F269").Value = 60
$ws.Range("G269").Value = 6072

$ws.Range("B274").Value = 96708.45

$ws.Range("F284").Value = 29
$ws.Range("G284").Value = 4202.97

$ws.Range("F285").Value = 8
$ws.Range("G285").Value = 1054.64

$ws.Range("B294").Value = 63571
$ws.Range("F294").Value = 5
$ws.Range("G294").Value = 717.4

$ws.Range("B295").Value = 63531
$ws.Range("F295").Value = 80
$ws.Range("G295").Value = 11478.4

$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12

$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 147
$ws.Range("G300").Value = 7003.08

$ws.Range("B311").Value = 61605
$ws.Range("E311").Value = 133.78
$ws.Range("F311").Value = -13
$ws.Range("G311").Value = -1455.48

$ws.Range("B312").Value = 63563
$ws.Range("E312").Value = 119.04
$ws.Range("F312").Value = 2
$ws.Range("G312").Value = 223.92

$ws.Range("F323").Value = 22
$ws.Range("G323").Value = 2609.86

$ws.Range("F324").Value = 47
$ws.Range("G324").Value = 2779.11

$ws.Range("F328").Value = 950
$ws.Range("G328").Value = 19978.5

$ws.Range("F329").Value = 39
$ws.Range("G329").Value = 6279

$ws.Range("F334").Value = 366
$ws.Range("G334").Value = 55328.22

$ws.Range("B339").Value = 327187.05

$ws.Range("F355").Value = 101
$ws.Range("G355").Value = 3256.24

$ws.Range("B361").Value = 13422.56

$ws.Range("F366").Value = 32
$ws.Range("G366").Value = 4394.24

$ws.Range("F382").Value = 159
$ws.Range("G382").Value = 6838.59

$ws.Range("F393").Value = 128
$ws.Range("G393").Value = 26993.92

$ws.Range("B395").Value = 253740.45

$ws.Range("F424").Value = 81
$ws.Range("G424").Value = 2408.94

$ws.Range("F426").Value = 132
$ws.Range("G426").Value = 12751.2

$ws.Range("B430").Value = 50686.39

$ws.Range("F470").Value = 171
$ws.Range("G470").Value = 2190.51

$ws.Range("F471").Value = 76
$ws.Range("G471").Value = 1998.8

$ws.Range("F484").Value = 607
$ws.Range("G484").Value = 3939.43

$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 214
$ws.Range("G485").Value = 2814.1

$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945

$ws.Range("F488").Value = 112
$ws.Range("G488").Value = 2945.6

$ws.Range("B490").Value = 53595
$ws.Range("E490").Value = 17.61
$ws.Range("F490").Value = -335
$ws.Range("G490").Value = -4934.55

$ws.Range("B491").Value = 65067
$ws.Range("E491").Value = 15.65
$ws.Range("F491").Value = 249
$ws.Range("G491").Value = 3667.77

$ws.Range("B492").Value = -243.38

$ws.Range("F516").Value = 18
$ws.Range("G516").Value = 291.24

$ws.Range("B528").Value = 18548.03

$ws.Range("F541").Value = 28
$ws.Range("G541").Value = 3087.56

$ws.Range("B542").Value = 9028.82

$ws.Range("F575").Value = 64
$ws.Range("G575").Value = 1696.64

$ws.Range("B582").Value = 22244.62

$ws.Range("F611").Value = 63
$ws.Range("G611").Value = 11054.61

$ws.Range("B614").Value = 42449.86

$ws.Range("F648").Value = 85
$ws.Range("G648").Value = 8834.9

$ws.Range("B651").Value = 38100.91

$ws.Range("F710").Value = 90
$ws.Range("G710").Value = 6260.4

$ws.Range("F713").Value = 381
$ws.Range("G713").Value = 51438.81

$ws.Range("F714").Value = 24
$ws.Range("G714").Value = 898.08

$ws.Range("F715").Value = 300
$ws.Range("G715").Value = 36213

$ws.Range("B716").Value = 159103.49

$ws.Range("F718").Value = 18
$ws.Range("G718").Value = 4490.64

$ws.Range("F721").Value = 8
$ws.Range("G721").Value = 870.48

$ws.Range("B732").Value = 65362
$ws.Range("F732").Value = 62
$ws.Range("G732").Value = 2533.94

$ws.Range("B733").Value = 65079
$ws.Range("F733").Value = 21
$ws.Range("G733").Value = 858.27

$ws.Range("F741").Value = 128
$ws.Range("G741").Value = 30950.4

$ws.Range("B743").Value = 90144.5

$ws.Range("F768").Value = 3260
$ws.Range("G768").Value = 531738.6

$ws.Range("F769").Value = 67
$ws.Range("G769").Value = 11803.39

$ws.Range("F771").Value = 495
$ws.Range("G771").Value = 71601.75

$ws.Range("F774").Value = 202
$ws.Range("G774").Value = 25973.16

$ws.Range("B775").Value = 813685.08

$ws.Range("F778").Value = 118
$ws.Range("G778").Value = 17229.18

$ws.Range("F789").Value = 41
$ws.Range("G789").Value = 1596.13

$ws.Range("F791").Value = 63
$ws.Range("G791").Value = 2490.39

$ws.Range("B792").Value = 82830.34

$ws.Range("B793").Value = 3102491.92

$ws.Range("B794").Value = 3102491.92
